$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(285).Insert()
$ws.Rows.Item(285).Insert()

$ws.Range("A285").Value = 9
$ws.Range("B285").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C285").Value = "Metropolitana"
$ws.Range("D285").Value = 44637
$ws.Range("E285").Value = 13
$ws.Range("F285").Value = 100112052
$ws.Range("G285").Value = "Albahaca"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 200
$ws.Range("K285").Value = 2500
$ws.Range("L285").Value = 3000
$ws.Range("M285").Value = 2750
$ws.Range("N285").Value = "`$/docena de matas"
$ws.Range("O285").Value = "Provincia de Cachapoal"
$ws.Range("P285").Value = 458
$ws.Range("Q285").Value = 6
$ws.Range("R285").Value = "Hortaliza"

$ws.Range("A286").Value = 9
$ws.Range("B286").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C286").Value = "Metropolitana"
$ws.Range("D286").Value = 44637
$ws.Range("E286").Value = 13
$ws.Range("F286").Value = 100112052
$ws.Range("G286").Value = "Albahaca"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 280
$ws.Range("K286").Value = 2500
$ws.Range("L286").Value = 3000
$ws.Range("M286").Value = 2679
$ws.Range("N286").Value = "`$/docena de matas"
$ws.Range("O286").Value = "Provincia de Chacabuco"
$ws.Range("P286").Value = 446
$ws.Range("Q286").Value = 6
$ws.Range("R286").Value = "Hortaliza"

Write-Host "done"
